$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wave3")

$ws.Range("A1").Value = "Baidu"
$ws.Range("A2").Select()
